$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 211-214 follow the exact same cell formatting as row 210
# (A/B/C/D center+bottom-ish font0 style, E font5 style, F default style).
# Copy that formatting down first, then fill in values.
$ws.Range("A210:F210").Copy()
$ws.Range("A211:F214").PasteSpecial(-4122)

# Row 211 - stream_search_phrase
$ws.Cells.Item(211, 1).Value = "stream_search_phrase"
$ws.Cells.Item(211, 2).Value = "src"
$ws.Cells.Item(211, 3).Value = "ui_utils_public"
$ws.Cells.Item(211, 4).Value = 901
$ws.Cells.Item(211, 5).Value = "Wybierz lub podaj frazę do wyszukania w wybranych stronach (ostatnie {num_hours}h)"
$ws.Cells.Item(211, 6).Value = "Select or enter a phrase to search in the chosen pages (last {num_hours}h)"

# Row 212 - stream_search_to_short_phrase
$ws.Cells.Item(212, 1).Value = "stream_search_to_short_phrase"
$ws.Cells.Item(212, 2).Value = "src"
$ws.Cells.Item(212, 3).Value = "ui_utils_public"
$ws.Cells.Item(212, 4).Value = 913
$ws.Cells.Item(212, 5).Value = "Zapytanie musi posiadać co majmniej {MIN_STREAM_QUERY_LEN} znaków"
$ws.Cells.Item(212, 6).Value = "The query must contain at least {MIN_STREAM_QUERY_LEN} characters"

# Row 213 - stream_searching
$ws.Cells.Item(213, 1).Value = "stream_searching"
$ws.Cells.Item(213, 2).Value = "src"
$ws.Cells.Item(213, 3).Value = "ui_utils_public"
$ws.Cells.Item(213, 4).Value = 918
$ws.Cells.Item(213, 5).Value = "Wyszukiwanie…"
$ws.Cells.Item(213, 6).Value = "Searching…"

# Row 214 - stream_searching_problem
$ws.Cells.Item(214, 1).Value = "stream_searching_problem"
$ws.Cells.Item(214, 2).Value = "src"
$ws.Cells.Item(214, 3).Value = "ui_utils_public"
$ws.Cells.Item(214, 4).Value = 931
$ws.Cells.Item(214, 5).Value = "Problem z poiłączeniem, spróbuj jeszcze raz"
$ws.Cells.Item(214, 6).Value = "Connection problem, please try again"

# Column A of rows 212-214 gets its own distinct style: a dedicated
# "Aptos Narrow" 11pt black font, centered horizontally and vertically.
$aRange = $ws.Range("A212:A214")
$aRange.HorizontalAlignment = -4108
$aRange.VerticalAlignment = -4108
$aRange.WrapText = $false
$aRange.Font.Name = "Aptos Narrow"
$aRange.Font.Size = 11
$aRange.Font.ColorIndex = 1
$aRange.Font.Bold = $false
$aRange.Font.Italic = $false
$aRange.Font.Underline = $false

# Scroll/selection state as left by the editor.
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("B213:C214").Select()

Write-Output "done"
